$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove existing hyperlinks before restructuring rows, to avoid stale/duplicate refs
$ws.Hyperlinks.Delete()

# Delete the last data row (row 8); rows 2-7 keep their row numbers
$ws.Rows.Item(8).Delete()

# Update changed cell values for rows 2-7 (new scrape batch at 06:28:57)
$ws.Range("A2").Value = "2025-11-04 06:28:57"
$ws.Range("B2").Value = "【急募】WordPressで施設検索サイトのMVPを相談しながら構築いただける方を探しています!"
$ws.Range("D2").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5426483"
$ws.Range("G2").Value = 65
$ws.Range("H2").Value = "◇サイト ○WordPress"
$ws.Range("A3").Value = "2025-11-04 06:28:57"
$ws.Range("B3").Value = "PHP業務アプリケーションの改修対応"
$ws.Range("D3").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5426598"
$ws.Range("G3").Value = 58
$ws.Range("H3").Value = "◇アプリ ○PHP"
$ws.Range("A4").Value = "2025-11-04 06:28:57"
$ws.Range("B4").Value = "社外エンジニア(WEBサイトやシステムのメンテナンス等の保守/改修等)の募集"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5426251"
$ws.Range("G4").Value = 53
$ws.Range("H4").Value = "◇サイト"
$ws.Range("A5").Value = "2025-11-04 06:28:57"
$ws.Range("B5").Value = "【急募】警備スタッフと各作業現場のマッチングシステム構築依頼"
$ws.Range("D5").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5426527"
$ws.Range("G5").Value = 40
$ws.Range("A6").Value = "2025-11-04 06:28:57"
$ws.Range("B6").Value = "〖リモート可〗Delphiエンジニア募集"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5341051"
$ws.Range("G6").Value = 25
$ws.Range("A7").Value = "2025-11-04 06:28:57"
$ws.Range("B7").Value = "PowerAutomate GoogleドライブからCSVをダウンロードしてヤマトWEBで印刷処理"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5426627"
$ws.Range("G7").Value = 13

# H5 has no skill-summary value in the updated data; clear any leftover content
$ws.Range("H5").ClearContents()

# Re-add hyperlinks for the URL column (F2:F7) pointing at their row URLs.
# (Use the same literal URL just written above rather than reading back
#  Range.Value, whose COM getter is unreliable in this host.)
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5426483") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5426598") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5426251") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5426527") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5341051") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5426627") | Out-Null

# Hyperlinks.Add() applies a freshly-duplicated style; reapply the built-in
# "Hyperlink" cell style so F2:F7 keep the original style index/formatting
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"

# Restore column widths: D (col 4) widened to 32, H (col 8) narrowed to 17
# (ColumnWidth setter adds a constant ~0.8333 padding offset vs. the value Excel stores,
#  so we pre-subtract 5/6 to land exactly on the target stored widths.)
$ws.Columns.Item(4).ColumnWidth = 32 - 5/6
$ws.Columns.Item(8).ColumnWidth = 17 - 5/6

Write-Output "done"
